$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.603.28'
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("D3").Value = '3.762.31'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.98'
$ws.Range("E5").Value = '  -0.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.21'
$ws.Range("E6").Value = '  -1.61%  '
$ws.Range("D7").Value = '3.759.07'
$ws.Range("E7").Value = '  -0.76%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -1.01%  '
$ws.Range("E10").Value = '  -2.97%  '
$ws.Range("E11").Value = '  -1.52%  '
$ws.Range("E12").Value = '  -1.11%  '
$ws.Range("E13").Value = '  -7.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.11'
$ws.Range("E14").Value = '  -1.80%  '
$ws.Range("D15").Value = '4.395.73'
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("D16").Value = '3.769.75'
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("D17").Value = '68.553.23'
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.94'
$ws.Range("E18").Value = '  -4.40%  '
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("E20").Value = '  -3.09%  '
$ws.Range("E21").Value = '  +1.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '465.54'
$ws.Range("E22").Value = '  -0.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.697'
$ws.Range("E23").Value = '  -3.40%  '
$ws.Range("E24").Value = '  -1.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.21'
$ws.Range("E25").Value = '  +0.48%  '
$ws.Range("E26").Value = '  -2.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.97'
$ws.Range("E27").Value = '  -1.70%  '
$ws.Range("E28").Value = '  -3.80%  '
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").Value = '3.910.06'
$ws.Range("E30").Value = '  -0.62%  '
$ws.Range("E31").Value = '  -4.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.34'
$ws.Range("E32").Value = '  -3.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.01'
$ws.Range("E33").Value = '  -1.98%  '
$ws.Range("E34").Value = '  -3.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.21'
$ws.Range("E35").Value = '  -0.81%  '
$ws.Range("D37").Value = '3.714.96'
$ws.Range("E37").Value = '  -0.77%  '
$ws.Range("E38").Value = '  -3.74%  '
$ws.Range("E39").Value = '  -9.09%  '
$ws.Range("E40").Value = '  -1.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.45%  '
$ws.Range("E42").Value = '  -1.10%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.302'
$ws.Range("E45").Value = '  -3.74%  '
$ws.Range("B46").Value = 'Arweave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '43.90'
$ws.Range("E46").Value = '  +8.38%  '
$ws.Range("E47").Value = '  +2.51%  '
$ws.Range("E48").Value = '  -2.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.49'
$ws.Range("E49").Value = '  -2.46%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '145.18'
$ws.Range("E50").Value = '  +2.11%  '
$ws.Range("B51").Value = 'Bittensor'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '390.32'
$ws.Range("E51").Value = '  -2.88%  '
